# Applies the cryptos list price/volume refresh described in the commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, [string]$text) {
    # Force Excel to store the value as literal text (matching the
    # inlineStr cells in the workbook) instead of silently parsing
    # numeric-looking strings (e.g. "0.4740") into numbers and losing
    # formatting such as trailing/leading zeros.
    $range.Value = "'" + $text
}

# --- Rows where the coin, link, price and volume all changed (re-ranked / swapped rows) ---
$ws.Range("B11").Value = "WrappedEther"
$ws.Range("C11").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D11").Value = "1.832.60"
$ws.Range("E11").Value = "  -1.13%  "

$ws.Range("B12").Value = "TRON"
$ws.Range("C12").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
Set-TextValue $ws.Range("D12") "0.07443"
$ws.Range("E12").Value = "  +0.04%  "

$ws.Range("B47").Value = "EnergySwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue $ws.Range("D47") "8.587"
$ws.Range("E47").Value = "  -1.16%  "

$ws.Range("B48").Value = "Elrond"
$ws.Range("C48").Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
Set-TextValue $ws.Range("D48") "33.47"
$ws.Range("E48").Value = "  +0.58%  "

# --- Rows where only price (D) and volume (E) changed ---
$ws.Range("D2").Value = "30.602.72"
$ws.Range("E2").Value = "  +1.26%  "
$ws.Range("D3").Value = "1.854.92"
$ws.Range("E3").Value = "  +0.01%  "
Set-TextValue $ws.Range("D4") "0.9992"
$ws.Range("E4").Value = "  -0.16%  "
Set-TextValue $ws.Range("D5") "234.05"
$ws.Range("E5").Value = "  +0.01%  "
Set-TextValue $ws.Range("D6") "0.9999"
$ws.Range("E6").Value = "  -0.14%  "
Set-TextValue $ws.Range("D7") "0.4740"
$ws.Range("E7").Value = "  +0.73%  "
Set-TextValue $ws.Range("D8") "0.2747"
$ws.Range("E8").Value = "  +0.90%  "
Set-TextValue $ws.Range("D9") "0.06335"
$ws.Range("E9").Value = "  -0.89%  "
Set-TextValue $ws.Range("D10") "17.63"
$ws.Range("E10").Value = "  +8.31%  "
Set-TextValue $ws.Range("D13") "5.014"
$ws.Range("E13").Value = "  +1.04%  "
Set-TextValue $ws.Range("D14") "84.63"
$ws.Range("E14").Value = "  -0.75%  "
Set-TextValue $ws.Range("D15") "0.6254"
$ws.Range("E15").Value = "  -0.76%  "
$ws.Range("D16").Value = "30.549.39"
$ws.Range("E16").Value = "  +1.19%  "
Set-TextValue $ws.Range("D17") "244.04"
$ws.Range("E17").Value = "  +6.05%  "
Set-TextValue $ws.Range("D18") "1.0000"
$ws.Range("E18").Value = "  -0.10%  "
Set-TextValue $ws.Range("D19") "12.67"
$ws.Range("E19").Value = "  -0.01%  "
Set-TextValue $ws.Range("D20") "0.000007357"
$ws.Range("E20").Value = "  +0.07%  "
Set-TextValue $ws.Range("D21") "0.9993"
$ws.Range("E21").Value = "  -0.21%  "
Set-TextValue $ws.Range("D22") "4.936"
$ws.Range("E22").Value = "  -0.43%  "
Set-TextValue $ws.Range("D23") "5.959"
$ws.Range("E23").Value = "  +0.01%  "
Set-TextValue $ws.Range("D24") "9.222"
$ws.Range("E24").Value = "  -0.69%  "
Set-TextValue $ws.Range("D25") "161.99"
$ws.Range("E25").Value = "  -2.60%  "
Set-TextValue $ws.Range("D26") "18.01"
$ws.Range("E26").Value = "  +0.81%  "
Set-TextValue $ws.Range("D27") "1.886"
$ws.Range("E27").Value = "  +0.54%  "
Set-TextValue $ws.Range("D28") "0.1025"
$ws.Range("E28").Value = "  -1.57%  "
Set-TextValue $ws.Range("D29") "1.373"
$ws.Range("E29").Value = "  -1.28%  "
Set-TextValue $ws.Range("D30") "4.032"
$ws.Range("E30").Value = "  -2.40%  "
Set-TextValue $ws.Range("D31") "3.830"
$ws.Range("E31").Value = "  -1.70%  "
Set-TextValue $ws.Range("D32") "0.04852"
$ws.Range("E32").Value = "  -1.13%  "
Set-TextValue $ws.Range("D33") "1.139"
$ws.Range("E33").Value = "  -2.07%  "
Set-TextValue $ws.Range("D34") "0.7039"
$ws.Range("E34").Value = "  -2.10%  "
Set-TextValue $ws.Range("D35") "2.704"
$ws.Range("E35").Value = "  -0.01%  "
Set-TextValue $ws.Range("D36") "0.01911"
$ws.Range("E36").Value = "  +2.18%  "
Set-TextValue $ws.Range("D37") "2.686"
$ws.Range("E37").Value = "  +1.62%  "
Set-TextValue $ws.Range("D38") "0.8759"
$ws.Range("E38").Value = "  -3.72%  "
Set-TextValue $ws.Range("D39") "1.991"
$ws.Range("E39").Value = "  +1.62%  "
Set-TextValue $ws.Range("D40") "106.25"
$ws.Range("E40").Value = "  +0.91%  "
Set-TextValue $ws.Range("D41") "0.9997"
$ws.Range("E41").Value = "  -0.02%  "
Set-TextValue $ws.Range("D42") "0.4076"
$ws.Range("E42").Value = "  -0.29%  "
Set-TextValue $ws.Range("D43") "5.532"
$ws.Range("E43").Value = "  -0.28%  "
Set-TextValue $ws.Range("D44") "7.212"
$ws.Range("E44").Value = "  +1.85%  "
Set-TextValue $ws.Range("D45") "62.50"
$ws.Range("E45").Value = "  +3.11%  "
Set-TextValue $ws.Range("D46") "0.1213"
$ws.Range("E46").Value = "  +1.65%  "
Set-TextValue $ws.Range("D49") "0.05551"
$ws.Range("E49").Value = "  -0.49%  "
Set-TextValue $ws.Range("D50") "1.362"
$ws.Range("E50").Value = "  -2.42%  "
Set-TextValue $ws.Range("D51") "0.3673"
$ws.Range("E51").Value = "  -0.36%  "
